$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must be forced to remain text
# (otherwise Excel auto-converts "580.38" style values into numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "62.210.67"
$ws.Cells.Item(2, 5).Value = "  -1.38%  "
$ws.Cells.Item(3, 4).Value = "2.451.56"
$ws.Cells.Item(3, 5).Value = "  -0.08%  "
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "580.38"
$ws.Cells.Item(5, 5).Value = "  +1.33%  "
$ws.Cells.Item(6, 4).Value = "143.54"
$ws.Cells.Item(6, 5).Value = "  -1.83%  "
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 5).Value = "  -1.09%  "
$ws.Cells.Item(9, 4).Value = "2.448.46"
$ws.Cells.Item(9, 5).Value = "  -0.09%  "
$ws.Cells.Item(10, 5).Value = "  -3.58%  "
$ws.Cells.Item(11, 5).Value = "  +2.16%  "
$ws.Cells.Item(12, 5).Value = "  -0.97%  "
$ws.Cells.Item(13, 4).Value = "0.344"
$ws.Cells.Item(13, 5).Value = "  -3.17%  "
$ws.Cells.Item(14, 4).Value = "26.51"
$ws.Cells.Item(14, 5).Value = "  -1.71%  "
$ws.Cells.Item(15, 5).Value = "  -3.74%  "
$ws.Cells.Item(16, 4).Value = "2.809.38"
$ws.Cells.Item(16, 5).Value = "  -3.17%  "
$ws.Cells.Item(17, 4).Value = "62.187.49"
$ws.Cells.Item(17, 5).Value = "  -1.39%  "
$ws.Cells.Item(18, 4).Value = "2.429.65"
$ws.Cells.Item(18, 5).Value = "  -0.74%  "
$ws.Cells.Item(19, 4).Value = "10.87"
$ws.Cells.Item(20, 5).Value = "  -2.57%  "
$ws.Cells.Item(21, 4).Value = "329.45"
$ws.Cells.Item(21, 5).Value = "  +0.33%  "
$ws.Cells.Item(22, 5).Value = "  -2.70%  "
$ws.Cells.Item(23, 4).Value = "1.99"
$ws.Cells.Item(23, 5).Value = "  -3.73%  "
$ws.Cells.Item(24, 5).Value = "  -0.56%  "
$ws.Cells.Item(25, 4).Value = "65.93"
$ws.Cells.Item(25, 5).Value = "  +1.09%  "
$ws.Cells.Item(26, 4).Value = "9.36"
$ws.Cells.Item(26, 5).Value = "  +5.94%  "
$ws.Cells.Item(27, 4).Value = "619.29"
$ws.Cells.Item(27, 5).Value = "  +0.84%  "
$ws.Cells.Item(28, 4).Value = "0.0₃0958"
$ws.Cells.Item(28, 5).Value = "  -5.97%  "
$ws.Cells.Item(29, 4).Value = "2.541.37"
$ws.Cells.Item(29, 5).Value = "  -1.52%  "
$ws.Cells.Item(30, 5).Value = "  +0.02%  "
$ws.Cells.Item(31, 4).Value = "1.44"
$ws.Cells.Item(31, 5).Value = "  -4.56%  "
$ws.Cells.Item(32, 4).Value = "8.01"
$ws.Cells.Item(32, 5).Value = "  -2.48%  "
$ws.Cells.Item(34, 5).Value = "  -1.57%  "
$ws.Cells.Item(35, 5).Value = "  -5.38%  "
$ws.Cells.Item(36, 5).Value = "  +0.14%  "
$ws.Cells.Item(37, 5).Value = "  -5.97%  "
$ws.Cells.Item(38, 4).Value = "0.376"
$ws.Cells.Item(38, 5).Value = "  -0.57%  "
$ws.Cells.Item(39, 4).Value = "150.25"
$ws.Cells.Item(39, 5).Value = "  +2.29%  "
$ws.Cells.Item(40, 4).Value = "5.31"
$ws.Cells.Item(40, 5).Value = "  -1.56%  "
$ws.Cells.Item(41, 5).Value = "  -2.67%  "
$ws.Cells.Item(42, 5).Value = "  -2.53%  "
$ws.Cells.Item(43, 4).Value = "42.51"
$ws.Cells.Item(43, 5).Value = "  +1.77%  "
$ws.Cells.Item(44, 5).Value = "  -0.02%  "
$ws.Cells.Item(45, 4).Value = "2.47"
$ws.Cells.Item(45, 5).Value = "  -5.01%  "
$ws.Cells.Item(46, 4).Value = "143.27"
$ws.Cells.Item(46, 5).Value = "  -3.68%  "
$ws.Cells.Item(47, 5).Value = "  -3.44%  "
$ws.Cells.Item(48, 5).Value = "  +0.53%  "
$ws.Cells.Item(49, 4).Value = "0.0524"
$ws.Cells.Item(49, 5).Value = "  -1.65%  "
$ws.Cells.Item(50, 4).Value = "0.0₆0242"
$ws.Cells.Item(50, 5).Value = "  +12.23%  "
$ws.Cells.Item(51, 5).Value = "  -7.35%  "

# Restore default (unstyled) cell style now that the text format has been applied,
# so the cells keep style index 0 like the rest of the sheet.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
